# Natmi following Dr Hou advice
# Rebuild the Wnt5a-Fzd5 LR-pair table: add the "sCs" sending-cluster rows
# and refresh the recomputed statistics for the existing FAPs rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2; A = "FAPs"; B = "Wnt5a"; C = "Fzd5"; D = "ECs"; Nums = @(3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 8.259748666666667, 24.779246, 0.4118513756345933, 0.4118513756345933, 86.18422746017555, 775.6580471415799, 0.3999526590379378, 0.3999526590379379) }
    @{ Row = 3; A = "FAPs"; B = "Wnt5a"; C = "Fzd5"; D = "FAPs"; Nums = @(3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 5.709686666666666, 17.12906, 0.2846990148258541, 0.2846990148258541, 59.57626003708889, 536.1863403337999, 0.2764738319245217, 0.2764738319245219) }
    @{ Row = 4; A = "FAPs"; B = "Wnt5a"; C = "Fzd5"; D = "M2"; Nums = @(3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 4.176559, 12.529677, 0.2082535000745028, 0.2082535000745028, 43.57923290202334, 392.2130961182099, 0.2022368894128777, 0.2022368894128777) }
    @{ Row = 5; A = "FAPs"; B = "Wnt5a"; C = "Fzd5"; D = "sCs"; Nums = @(3, 1, 10.43424333333333, 31.30273, 0.9711091978791583, 0.9711091978791584, 3, 1, 1.909174, 5.727522, 0.09519610946504979, 0.09519610946504979, 19.92078608167333, 179.28707473506, 0.09244581750382105, 0.09244581750382107) }
    @{ Row = 6; A = "sCs"; B = "Wnt5a"; C = "Fzd5"; D = "ECs"; Nums = @(2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 8.259748666666667, 24.779246, 0.4118513756345933, 0.4118513756345933, 2.564007700604, 23.076069305436, 0.01189871659665544, 0.01189871659665544) }
    @{ Row = 7; A = "sCs"; B = "Wnt5a"; C = "Fzd5"; D = "FAPs"; Nums = @(2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 5.709686666666666, 17.12906, 0.2846990148258541, 0.2846990148258541, 1.77241235444, 15.95171118996, 0.008225182901332301, 0.008225182901332303) }
    @{ Row = 8; A = "sCs"; B = "Wnt5a"; C = "Fzd5"; D = "M2"; Nums = @(2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 4.176559, 12.529677, 0.2082535000745028, 0.2082535000745028, 1.296495797898, 11.668462181082, 0.006016610661625134, 0.006016610661625135) }
    @{ Row = 9; A = "sCs"; B = "Wnt5a"; C = "Fzd5"; D = "sCs"; Nums = @(2, 0.6666666666666666, 0.310422, 0.9312659999999999, 0.02889080212084161, 0.02889080212084161, 3, 1, 1.909174, 5.727522, 0.09519610946504979, 0.09519610946504979, 0.592649611428, 5.333846502852, 0.00275029196122873, 0.00275029196122873) }
)

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    for ($i = 0; $i -lt $numCols.Length; $i++) {
        $ws.Range("$($numCols[$i])$row").Value = $r.Nums[$i]
    }
}

